# Append the latest "逃离鸭科夫" mod-count entry to the ModCounts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row (row 8) and target the row right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Carry the centered-alignment style used by the existing data rows onto the
# new row.
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 3))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 3))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

# Build the literal date text off-sheet (forcing Text format there) so typing
# "2025/11/18" doesn't get auto-converted into a date serial number, then
# bring over just the computed value - leaving the destination cell's
# (already-copied) style untouched.
$scratch = $ws.Cells.Item($ws.Rows.Count, $ws.Columns.Count)
$scratch.NumberFormat = "@"
$scratch.Value = "2025/11/18"
$scratch.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1167
